# Appends the "7, 9, 11, ... 55" checker-count rows (rows 8-37) to the
# "guerrilla vs. hardcoded COIN 6 checkers" results sheet, extending the
# used range from A1:C7 to A1:C37.
#
# Column A holds the checker count as TEXT (to match the existing rows,
# which were written as text-typed cells, not numbers) styled like the
# other header/label cells in column A (bold, bordered, centered);
# columns B/C hold the win-rate and average-game-length as plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 is the last pre-existing data row; its A-cell already carries the
# exact label style (bold font, thin border, centered alignment) used
# throughout column A. Re-using it as a paste-format template guarantees
# every newly written A-cell ends up on that very same style, rather than
# a freshly synthesized (and therefore not-identical) one.
$template = $ws.Range("A7")

# row, checker-count label (text), win rate, avg. game length
$data = @(
    @(8, "7", 0, 33),
    @(9, "9", 0, 3),
    @(10, "11", 0, 33),
    @(11, "13", 0, 33),
    @(12, "15", 100, 23),
    @(13, "17", 0, 22),
    @(14, "19", 0, 33),
    @(15, "21", 0, 22),
    @(16, "23", 100, 28),
    @(17, "25", 0, 33),
    @(18, "27", 0, 33),
    @(19, "29", 0, 21),
    @(20, "31", 0, 23),
    @(21, "33", 0, 33),
    @(22, "35", 0, 33),
    @(23, "37", 0, 24),
    @(24, "39", 100, 31),
    @(25, "41", 0, 22),
    @(26, "43", 0, 33),
    @(27, "45", 0, 33),
    @(28, "46", 0, 1),
    @(29, "47", 0, 33),
    @(30, "48", 0, 1),
    @(31, "49", 0, 33),
    @(32, "50", 0, 1),
    @(33, "51", 100, 25),
    @(34, "52", 100, 31),
    @(35, "53", 0, 33),
    @(36, "54", 0, 1),
    @(37, "55", 0, 21)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $label = $entry[1]
    $winRate = $entry[2]
    $avgLength = $entry[3]

    # Force column A to be stored as text (General-formatted cells would
    # otherwise silently coerce a numeric-looking string back into a
    # number), then stamp the exact label style on top of it.
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $label
    $template.Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("B$r").Value = $winRate
    $ws.Range("C$r").Value = $avgLength
}
